# "Implementação do HighLight nos botões e documentação"
# Rows 55 ("Tratar highlight nos botões:"), 59 ("Pausa") and 61 ("Fim do jogo")
# get actual-duration (column B) entries and are marked as completed ("S",
# shared string idx 73) in column I instead of the inherited "N" formula
# result. This ripples through the SUM/SUMIF totals and EVM ratios further
# down the sheet (rows 64 and 68-71), which the engine recalculates for us.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

$ws.Range("B55").Value = 0.5
$ws.Range("I55").Value = "S"

$ws.Range("B59").Value = 2
$ws.Range("I59").Value = "S"

$ws.Range("B61").Value = 2
$ws.Range("I61").Value = "S"

# Match the author's new scroll position / active cell in the sheet view.
$win = $excel.ActiveWindow
$win.ScrollRow = 50
$win.ScrollColumn = 1
$ws.Range("I56").Select()
